$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the new calculated-column table (H1:L1)
$ws.Range("H1").Value = "LEFT"
$ws.Range("I1").Value = "RIGHT"
$ws.Range("J1").Value = "MID"
$ws.Range("K1").Value = "CONCAT"
$ws.Range("L1").Value = "TRIM"

# Turn H1:L31 into a second Excel table (Table2)
$lo2 = $ws.ListObjects.Add(1, $ws.Range("H1:L31"), 0, 1)
$lo2.Name = "Table2"
$lo2.TableStyle = "TableStyleMedium12"

# Fill in the calculated-column formulas for every data row (2-31).
# Cell-by-cell assignment keeps each formula as its own (non-shared) formula,
# matching how Excel stores calculated-table-column formulas.
for ($r = 2; $r -le 31; $r++) {
    $ws.Range("H$r").Formula = "=LEFT(Table1[[#This Row],[Product Codes]],5)"
    $ws.Range("I$r").Formula = "=RIGHT(Table1[[#This Row],[Product Codes]],4)"
    $ws.Range("J$r").Formula = "=MID(Table1[[#This Row],[Client Code]],4,2)"
    $ws.Range("K$r").Formula = "=CONCAT(Table2[[#This Row],[LEFT]],Table2[[#This Row],[RIGHT]])"
    $ws.Range("L$r").Formula = "=TRIM(Table1[[#This Row],[Client]])"
}

# Widen the new TRIM column and move the active selection onto it
$ws.Columns.Item(12).ColumnWidth = 20
$ws.Range("L3").Select()
